# Update the existing "Publisher" sheet data/headers and then add a new
# "Media" sheet with Journal/Publisher/Indexer data, matching the target
# workbook layout.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Publisher")

# --- Update header row (bold/yellow style s="2" is preserved automatically
# because we are only changing cell values, not formatting) ---
$ws1.Range("A1").Value = "Publisher Name"
$ws1.Range("B1").Value = "Place Of Publication"
$ws1.Range("C1").Value = "Country"

# --- Update data rows ---
$ws1.Range("A2").Value = "Manoj"
$ws1.Range("B2").Value = "Noida"
$ws1.Range("C2").Value = "China"

$ws1.Range("A3").Value = "Prashant"
$ws1.Range("B3").Value = "Christ the Redeemer"
$ws1.Range("C3").Value = "Brazil"

$ws1.Range("A4").Value = "Heena"
$ws1.Range("B4").Value = "Copacabana"
$ws1.Range("C4").Value = "France"

# Leave the selection on Publisher parked at C4 (matches the saved file,
# which no longer has Publisher as the visible/active tab).
$ws1.Range("C4").Select()

# --- Add the new "Media" worksheet right after "Publisher" ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Media"

$mediaHeaders = @("Journal Name", "Publisher Name", "Indexer")
$mediaRows = @(
  @("Test 11", "Pub1", "neha"),
  @("Test 22", "Pub",  "neha"),
  @("Test 3",  "Pub2", "neha"),
  @("Test 4",  "Pub3", "neha")
)

for ($c = 0; $c -lt $mediaHeaders.Length; $c++) {
  $ws2.Cells.Item(1, $c + 1).Value = $mediaHeaders[$c]
}
for ($r = 0; $r -lt $mediaRows.Length; $r++) {
  for ($c = 0; $c -lt $mediaRows[$r].Length; $c++) {
    $ws2.Cells.Item($r + 2, $c + 1).Value = $mediaRows[$r][$c]
  }
}

# Header row styling to match the bold + yellow fill used on Publisher's
# header row (reuses the existing style, no new style rows are created).
$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:C1").Interior.Color = 65535

# Autosize the first two columns like the Publisher sheet.
$ws2.Columns("A:B").AutoFit()

# Make Media the active/visible sheet and park its selection at E6, matching
# the saved workbook state.
$ws2.Activate()
$ws2.Range("E6").Select()
